$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 (columns A:AY) get cyclically permuted:
#   new row2 <- old row3
#   new row3 <- old row5
#   new row4 <- old row2
#   new row5 <- old row4
# This is implemented by reading each row's values and writing them back
# in the new arrangement, which correctly fills/clears columns that differ
# between rows (e.g. K and AC are present on some rows but not others).

# Columns I (counts like "4") and Y/AA (dates like "2023-08-11") hold text
# values that look numeric/date-like. Excel's automatic type inference on
# assignment would otherwise silently turn them into real numbers/dates, so
# force just those columns to a Text number format before writing the
# values back. (The other text columns round-trip fine as strings already.)
$textColumns = @("I","Y","AA")
foreach ($col in $textColumns) {
    $ws.Range($col + "2:" + $col + "5").NumberFormat = "@"
}

$old2 = $ws.Range("A2:AY2").Value2
$old3 = $ws.Range("A3:AY3").Value2
$old4 = $ws.Range("A4:AY4").Value2
$old5 = $ws.Range("A5:AY5").Value2

$ws.Range("A2:AY2").Value2 = $old3
$ws.Range("A3:AY3").Value2 = $old5
$ws.Range("A4:AY4").Value2 = $old2
$ws.Range("A5:AY5").Value2 = $old4

# Row 3 now holds what used to be row 5's data, which had no K/AC values.
# Writing the captured blank (Empty) values above turns those destination
# cells into present-but-empty cells (since they still carry the Text
# NumberFormat applied earlier). Use Clear() (contents + formatting) so the
# cells are fully removed from the sheet, matching the source row's shape.
$ws.Range("K3").Clear()
$ws.Range("AC3").Clear()
